# Apply updates to Jogos_da_Semana_FlashScore_2024-11-14.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 update
$ws.Range("W9").Value = 5.5

# Row 10 updates
$ws.Range("J10").Value = 2.37
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("Q10").Value = 2.2
$ws.Range("R10").Value = 1.65
$ws.Range("AI10").Value = 26
$ws.Range("AK10").Value = 51
$ws.Range("AL10").Value = 41
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 29
